$d = $word.ActiveDocument

function Split-ParagraphTextIntoRuns($paragraph, [string[]]$tokens) {
    $r = $paragraph.Range
    # Exclude the trailing paragraph mark from the range we replace.
    $target = $d.Range($r.Start, $r.End - 1)

    $runsXml = ""
    foreach ($tok in $tokens) {
        $escaped = $tok.Replace("&", "&amp;").Replace("<", "&lt;").Replace(">", "&gt;")
        $runsXml += '<w:r><w:t xml:space="preserve">' + $escaped + '</w:t></w:r>'
    }

    $xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' + `
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
        '<pkg:xmlData>' + `
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
        '<w:body><w:p>' + $runsXml + '</w:p></w:body>' + `
        '</w:document>' + `
        '</pkg:xmlData></pkg:part></pkg:package>'

    $target.InsertXML($xml)
}

# --- Title paragraph ---
$titleTokens = @("Questions:", " ", "Trigonometric", " ", "identities", " ", "(degrees)")
Split-ParagraphTextIntoRuns $d.Paragraphs.Item(1) $titleTokens

# --- Author paragraph ---
$authorTokens = @("Dzhemma", " ", "Ruseva")
Split-ParagraphTextIntoRuns $d.Paragraphs.Item(2) $authorTokens

# --- Abstract paragraph ---
$abstractTokens = @(
    "A", " ", "selection", " ", "of", " ", "questions", " ", "on", " ",
    "trigonometric", " ", "identities,", " ", "where", " ", "angles", " ",
    "are", " ", "measured", " ", "in", " ", "degrees."
)
Split-ParagraphTextIntoRuns $d.Paragraphs.Item(4) $abstractTokens

Write-Host "Edit complete"
